# Insert a new data row at row 423 (pushes existing rows 423:511 down to 424:512)
# and populate it with the new record for "Región de Arica y Parinacota".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A423").EntireRow.Insert()

$ws.Cells.Item(423, 1).Value = 3
$ws.Cells.Item(423, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(423, 3).Value = "Coquimbo"
$ws.Cells.Item(423, 4).Value = 44782
$ws.Cells.Item(423, 5).Value = 5
$ws.Cells.Item(423, 6).Value = 100112032
$ws.Cells.Item(423, 7).Value = "Zapallo italiano"
$ws.Cells.Item(423, 8).Value = "Sin especificar"
$ws.Cells.Item(423, 9).Value = "Primera"
$ws.Cells.Item(423, 10).Value = 165
$ws.Cells.Item(423, 11).Value = 16000
$ws.Cells.Item(423, 12).Value = 17000
$ws.Cells.Item(423, 13).Value = 16485
$ws.Cells.Item(423, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(423, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(423, 16).Value = 236
$ws.Cells.Item(423, 17).Value = 70
$ws.Cells.Item(423, 18).Value = "Hortaliza"

# Match the date number format used by the other rows in column D
$ws.Cells.Item(423, 4).NumberFormat = $ws.Cells.Item(424, 4).NumberFormat
